$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns (copy formatting from H1, then set text)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I0 (col I) and IF (col J), rows 2-43
$data = @(
    @(2, 8, 8),
    @(3, 9, 9),
    @(4, 7, 8),
    @(5, 7, 8),
    @(6, 8, 8),
    @(7, 4, 7),
    @(8, 9, 9),
    @(9, 13, 13),
    @(10, 7, 7),
    @(11, 6, 6),
    @(12, 7, 8),
    @(13, 9, 9),
    @(14, 8, 8),
    @(15, 8, 8),
    @(16, 8, 8),
    @(17, 7, 7),
    @(18, 8, 9),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 7, 8),
    @(23, 8, 8),
    @(24, 9, 9),
    @(25, 7, 7),
    @(26, 7, 8),
    @(27, 6, 6),
    @(28, 6, 7),
    @(29, 7, 7),
    @(30, 9, 9),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 7, 8),
    @(34, 7, 8),
    @(35, 8, 8),
    @(36, 5, 6),
    @(37, 7, 7),
    @(38, 7, 8),
    @(39, 9, 9),
    @(40, 7, 8),
    @(41, 4, 4),
    @(42, 3, 3),
    @(43, 5, 5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $i0 = $entry[1]
    $if_ = $entry[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if_
}
